$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph via Find.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Ver no Jupiter...' paragraph to remove."
}

# Work out which paragraph (by index) the found text belongs to, so we can
# also reach its immediate neighbours.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not resolve paragraph index for the found text."
}

# Remove three whole paragraphs:
#   - the blank paragraph right before "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph itself
#   - the following "© 2020 . Contact: ..." paragraph
# (The blank paragraph that used to follow the copyright line remains,
#  immediately preceding the page-break paragraph.)
$startPara = $d.Paragraphs.Item($targetIndex - 1)
$endPara = $d.Paragraphs.Item($targetIndex + 1)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()
